$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").ClearContents()
$ws.Range("I2").ClearContents()

# Row 3
$ws.Range("D3").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("K3").ClearContents()

# Row 4
$ws.Range("D4").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("K4").ClearContents()

# Row 5
$ws.Range("D5").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("K5").ClearContents()

# Row 6
$ws.Range("D6").ClearContents()
$ws.Range("H6").ClearContents()
$ws.Range("K6").ClearContents()

# Row 7
$ws.Range("D7").ClearContents()

# Row 8
$ws.Range("D8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("K8").ClearContents()

# Row 9
$ws.Range("D9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("K9").ClearContents()

# Row 10
$ws.Range("D10").ClearContents()
$ws.Range("H10").ClearContents()
$ws.Range("K10").ClearContents()

# Row 11
$ws.Range("D11").ClearContents()
$ws.Range("H11").ClearContents()
$ws.Range("K11").ClearContents()

# Row 12
$ws.Range("D12").ClearContents()
$ws.Range("K12").ClearContents()

# Row 13
$ws.Range("D13").ClearContents()
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 851,5 TL"
$ws.Range("H13").ClearContents()
$ws.Range("K13").ClearContents()

# Row 14
$ws.Range("D14").ClearContents()
$ws.Range("H14").ClearContents()
$ws.Range("K14").ClearContents()

# Row 24
$ws.Range("F24").ClearContents()

# Row 25
$ws.Range("F25").ClearContents()
